$d = $word.ActiveDocument

# Locate the "Ver no Jupiter..." paragraph and the preceding empty paragraph,
# plus the "(c) 2020 ..." paragraph that follows it, and remove all three,
# collapsing them away while leaving the surrounding paragraphs untouched.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $i
        break
    }
}

if ($null -ne $target) {
    $startPara = $d.Paragraphs.Item($target - 1)   # the blank paragraph right before it
    $endPara   = $d.Paragraphs.Item($target + 1)    # the "(c) 2020 ..." paragraph right after it

    $start = $startPara.Range.Start
    $end   = $endPara.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
